$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 656, shifting the
# existing weekly records (old rows 656-749) down to rows 658-751.
$ws.Rows.Item(656).Insert()
$ws.Rows.Item(656).Insert()

# Seed the two new rows with the same "shape" as the records that now sit
# directly below them (same product / grade / unit / region metadata), then
# overwrite the weekly figures (date, min price, average price, dozen price)
# with the new week's values.
$ws.Range("A658:R658").Copy()
$ws.Range("A656").PasteSpecial()

$ws.Range("A659:R659").Copy()
$ws.Range("A657").PasteSpecial()

$ws.Range("D656").Value = 44776
$ws.Range("J656").Value = 2900
$ws.Range("M656").Value = 6586
$ws.Range("P656").Value = 1098

$ws.Range("D657").Value = 44776
$ws.Range("J657").Value = 1550
$ws.Range("M657").Value = 4774
$ws.Range("P657").Value = 796
